# "Aggiunti alcuni grafici + consuntivo"
# Consuntivo: for the roles that had no hours/cost recorded ("-"),
# enter an actual numeric cost value of 0 instead of the placeholder text.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Riga 4 (Progettista) e Riga 6 (Programmatore): costo consuntivato a 0
$ws.Range("C4").Value = 0
$ws.Range("C6").Value = 0

# Selezione finale sull'intera tabella (come da salvataggio originale)
$ws.Range("A1:C8").Select()
